$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 185
$ws.Range("I19").Value = 180
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 180
$ws.Range("L19").Value = 200
$ws.Range("M19").Value = -5
$ws.Range("N19").Value = -550
$ws.Range("H40").Value = 3405.2632
$ws.Range("J40").Value = 3823.077
$ws.Range("L40").Value = 3823.077
$ws.Range("N40").Value = -4173.077
$ws.Range("H74").Value = 3812.75
$ws.Range("I74").Value = 3125.5
$ws.Range("K74").Value = 3125.5
$ws.Range("M74").Value = -2189.5
$ws.Range("H77").Value = 3812.75
$ws.Range("I77").Value = 3125.5
$ws.Range("K77").Value = 15627.5
$ws.Range("M77").Value = -10947.5
$ws.Range("H98").Value = 3546.44
$ws.Range("I98").Value = 2615.611
$ws.Range("J98").Value = 5940
$ws.Range("K98").Value = 2615.611
$ws.Range("L98").Value = 5940
$ws.Range("M98").Value = -1117.611
$ws.Range("N98").Value = -8936
$ws.Range("H122").Value = 3546.44
$ws.Range("I122").Value = 2615.611
$ws.Range("J122").Value = 5940
$ws.Range("K122").Value = 7846.833
$ws.Range("L122").Value = 17820
$ws.Range("M122").Value = -5396.833
$ws.Range("N122").Value = -22720
$ws.Range("H129").Value = 1620.25
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 2628.4038
$ws.Range("J132").Value = 5833
$ws.Range("L132").Value = 17499
$ws.Range("N132").Value = -22559
$ws.Range("H137").Value = 3649.647
$ws.Range("I137").Value = 2064.818
$ws.Range("K137").Value = 6194.454000000001
$ws.Range("M137").Value = -3644.454000000001
$ws.Range("H138").Value = 3659.9683
$ws.Range("I138").Value = 3542.158
$ws.Range("J138").Value = 3710.8408
$ws.Range("K138").Value = 10626.474
$ws.Range("L138").Value = 11132.5224
$ws.Range("M138").Value = -5486.474
$ws.Range("N138").Value = -21412.5224
$ws.Range("H141").Value = 15529.5
$ws.Range("I141").Value = 22265
$ws.Range("J141").Value = 5426.25
$ws.Range("K141").Value = 66795
$ws.Range("L141").Value = 16278.75
$ws.Range("M141").Value = -61615
$ws.Range("N141").Value = -26638.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16562.62
$ws.Range("I32").Value = 12062.709
$ws.Range("J32").Value = 47499.5
$ws.Range("K32").Value = 12062.709
$ws.Range("L32").Value = 47499.5
$ws.Range("M32").Value = -11775.709
$ws.Range("N32").Value = -48073.5
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H122").Value = 4998.4443
$ws.Range("I122").Value = 4992
$ws.Range("J122").Value = 4999.25
$ws.Range("K122").Value = 14976
$ws.Range("L122").Value = 14997.75
$ws.Range("M122").Value = -12526
$ws.Range("N122").Value = -19897.75
$ws.Range("H135").Value = 91582.836
$ws.Range("J135").Value = 91582.836
$ws.Range("L135").Value = 91582.836
$ws.Range("N135").Value = -101722.836

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 19187.2
$ws.Range("I97").Value = 19187.2
$ws.Range("K97").Value = 19187.2
$ws.Range("M97").Value = -18196.2
$ws.Range("H139").Value = 74897.5
$ws.Range("J139").Value = 99796
$ws.Range("L139").Value = 99796
$ws.Range("N139").Value = -110076

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3739.1738
$ws.Range("I31").Value = 2943.8125
$ws.Range("J31").Value = 5557.143
$ws.Range("K31").Value = 2943.8125
$ws.Range("L31").Value = 5557.143
$ws.Range("M31").Value = -2648.8125
$ws.Range("N31").Value = -6147.143
$ws.Range("H34").Value = 3739.1738
$ws.Range("I34").Value = 2943.8125
$ws.Range("J34").Value = 5557.143
$ws.Range("K34").Value = 2943.8125
$ws.Range("L34").Value = 5557.143
$ws.Range("M34").Value = -2741.8125
$ws.Range("N34").Value = -5961.143
$ws.Range("H132").Value = 266263.72
$ws.Range("I132").Value = 3106.1667
$ws.Range("K132").Value = 9318.500100000001
$ws.Range("M132").Value = -6788.500100000001
$ws.Range("H134").Value = 4681.6743
$ws.Range("I134").Value = 3925.6875
$ws.Range("J134").Value = 6880.909
$ws.Range("K134").Value = 11777.0625
$ws.Range("L134").Value = 20642.727
$ws.Range("M134").Value = -9242.0625
$ws.Range("N134").Value = -25712.727

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 833.38464
$ws.Range("I2").Value = 861.1667
$ws.Range("K2").Value = 5167.0002
$ws.Range("M2").Value = -5054.0002
$ws.Range("H34").Value = 1940.8125
$ws.Range("J34").Value = 2210.6428
$ws.Range("L34").Value = 6631.928400000001
$ws.Range("N34").Value = -6799.928400000001
$ws.Range("H55").Value = 1900.8
$ws.Range("J55").Value = 3000
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9354
$ws.Range("H56").Value = 5600
$ws.Range("I56").Value = 5600
$ws.Range("K56").Value = 5600
$ws.Range("M56").Value = -5070
$ws.Range("H132").Value = 2603.782
$ws.Range("I132").Value = 1420.3334
$ws.Range("J132").Value = 2758.145
$ws.Range("K132").Value = 12783.0006
$ws.Range("L132").Value = 24823.305
$ws.Range("M132").Value = -10253.0006
$ws.Range("N132").Value = -29883.305

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 656.3333
$ws.Range("I2").Value = 540.4
$ws.Range("K2").Value = 540.4
$ws.Range("M2").Value = -427.4
$ws.Range("H97").Value = 1420.4
$ws.Range("I97").Value = 1531.625
$ws.Range("K97").Value = 1531.625
$ws.Range("M97").Value = -1035.625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 767.1667
$ws.Range("I22").Value = 665.4286
$ws.Range("K22").Value = 665.4286
$ws.Range("M22").Value = -370.4286
$ws.Range("H27").Value = 767.1667
$ws.Range("I27").Value = 665.4286
$ws.Range("K27").Value = 665.4286
$ws.Range("M27").Value = -558.4286
$ws.Range("H46").Value = 306549.6
$ws.Range("I46").Value = 2920.7144
$ws.Range("J46").Value = 388295.84
$ws.Range("K46").Value = 2920.7144
$ws.Range("L46").Value = 388295.84
$ws.Range("M46").Value = -2732.7144
$ws.Range("N46").Value = -388671.84
$ws.Range("H55").Value = 316.66666
$ws.Range("I55").Value = 320
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 320
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = -147
$ws.Range("N55").Value = -646
$ws.Range("H132").Value = 134745.19
$ws.Range("I132").Value = 253432.38
$ws.Range("J132").Value = 6434.7295
$ws.Range("K132").Value = 760297.14
$ws.Range("L132").Value = 19304.1885
$ws.Range("M132").Value = -757767.14
$ws.Range("N132").Value = -24364.1885
$ws.Range("H136").Value = 6820.511
$ws.Range("I136").Value = 7030.4595
$ws.Range("J136").Value = 5849.5
$ws.Range("K136").Value = 21091.3785
$ws.Range("L136").Value = 17548.5
$ws.Range("M136").Value = -18541.3785
$ws.Range("N136").Value = -22648.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 76996.336
$ws.Range("J27").Value = 76996.336
$ws.Range("L27").Value = 76996.336
$ws.Range("N27").Value = -77134.336
$ws.Range("H97").Value = 25286.5
$ws.Range("J97").Value = 25286.5
$ws.Range("L97").Value = 25286.5
$ws.Range("N97").Value = -27268.5
$ws.Range("H132").Value = 212687.55
$ws.Range("I132").Value = 248217.88
$ws.Range("J132").Value = 4581.2856
$ws.Range("K132").Value = 744653.64
$ws.Range("L132").Value = 13743.8568
$ws.Range("M132").Value = -742123.64
$ws.Range("N132").Value = -18803.8568
$ws.Range("H133").Value = 121404.58
$ws.Range("J133").Value = 125927.055
$ws.Range("L133").Value = 125927.055
$ws.Range("N133").Value = -136047.055
$ws.Range("H136").Value = 9018.233
$ws.Range("I136").Value = 9321.879999999999
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 27965.64
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -25415.64
$ws.Range("N136").Value = -27600
